$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("boson"), shifting existing
# columns E..S to F..T.
$ws.Columns("E:E").Insert() | Out-Null

# New column header (row 1) and values (rows 2-12) for "pt_max".
$ws.Range("E1").Value = "pt_max"
$ws.Range("E2:E12").Value = 50

# Restore the selection to match the user's final selection (E2:E12).
$ws.Range("E2:E12").Select() | Out-Null
